# Fixed bug #31: Error if no file is imported on the bootstrap. Updated BM as well.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug Log")

# --- Update row 31 (bug #29 "Drop Bid round 2") - mark as Resolved ---
$ws.Range("E30").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = "Resolved"
$ws.Range("G31").Value = "14/11/2019"
$ws.Range("H30").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("H31").Value = "Brendon & Vittorio"

# --- Row 32: new bug #30 "Admin Page" (only S/N, Iteration, Function, Description filled) ---
$ws.Range("A30:D30").Copy()
$ws.Range("A32:D32").PasteSpecial(-4122)
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = 4
$ws.Range("C32").Value = "Admin Page"
$ws.Range("D32").Value = "Exclamation mark appears if no file is imported on the admin page"

# --- Row 33: new bug #31 "Admin Page" / bootstrap fix (fully filled) ---
$ws.Range("A29:E29").Copy()
$ws.Range("A33:E33").PasteSpecial(-4122)
$ws.Range("F31").Copy()
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("F31").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("H29").Copy()
$ws.Range("H33").PasteSpecial(-4122)

$ws.Range("A33").Value = 31
$ws.Range("B33").Value = 4
$ws.Range("C33").Value = "Admin Page"
$ws.Range("D33").Value = "Error if no file is imported on the bootstrap "
$ws.Range("E33").Value = "Resolved"
$ws.Range("F33").Value = "16/11/2019"
$ws.Range("G33").Value = "16/11/2019"
$ws.Range("H33").Value = "Matthew & DaEun"

$ws.Range("A33").Select()
